# Attempt to do table expansion again
#
# Adds a new "Replace table" config block to the Config sheet (rows 34-37):
#   Name         is  Replace table     (Copy from one table to another)
#   Table        is  RangleTable       ((Yes, it's a typo))
#   Target table is  PlanTable
#   Expand       is  TRUE
#
# and updates the sheet selections: Summary becomes the active/selected tab
# (cell A1), while Config's selection moves down to the newly added block.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Summary")
$ws2 = $wb.Worksheets.Item("Config")

# --- New config rows 34-37 on the Config sheet ---
$ws2.Range("B34").Value = "Name"
$ws2.Range("C34").Value = "is"
$ws2.Range("D34").Value = "Replace table"
$ws2.Range("F34").Value = "Copy from one table to another"

$ws2.Range("B35").Value = "Table"
$ws2.Range("C35").Value = "is"
$ws2.Range("D35").Value = "RangleTable"
$ws2.Range("F35").Value = "(Yes, it’s a typo)"

$ws2.Range("B36").Value = "Target table"
$ws2.Range("C36").Value = "is"
$ws2.Range("D36").Value = "PlanTable"

$ws2.Range("B37").Value = "Expand"
$ws2.Range("C37").Value = "is"
$ws2.Range("D37").Value = $true

# Column B holds the bold "key" labels, matching the rest of the sheet.
$ws2.Range("B34:B37").Font.Bold = $true

# --- Selection / active sheet changes ---
$ws2.Range("B40").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("A1").Select() | Out-Null
